# "Generate Report for Handoff"
#
# This report workbook has three sheets:
#   1) Overview  - one row per source file, with a "Status" column duplicated
#                  for each target locale (zh-cn in col E, de-de in col F)
#                  and a "Latest HO Xliff Generate Date" column (col G).
#   2) zh-cn     - per-locale detail sheet; Status is col C, "Latest Handoff
#                  Datetime" is col H.
#   3) de-de     - same layout as zh-cn.
#
# The handoff-generation run flips the status text from the old "handed
# back" message to "Ready for handoff" and refreshes the associated
# timestamps.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text + refreshed timestamps -----------------------------------

# Overview sheet: status repeated for each locale column, plus the shared
# "Latest HO Xliff Generate Date" timestamp.
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-29 12:59:47"

# zh-cn detail sheet: status + its "Latest Handoff Datetime".
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-29 12:59:42"

# de-de detail sheet: status + its "Latest Handoff Datetime".
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-29 12:59:47"

# --- Column width adjustments ----------------------------------------------
# The Status columns are narrowed now that the text is shorter
# ("Ready for handoff" vs. "Handed back: in sync with en-US").
# Target stored column width (OOXML character units) is 17.2159881591797.
# This host quantizes ColumnWidth assignments to whole pixels (steps of
# 1/6 of a character unit), so we feed it the input value that rounds to
# the closest representable width to the target.
$targetColWidth = 16.333333333333336

$wsOverview.Columns.Item("E").ColumnWidth = $targetColWidth
$wsOverview.Columns.Item("F").ColumnWidth = $targetColWidth
$wsZhCn.Columns.Item("C").ColumnWidth = $targetColWidth
$wsDeDe.Columns.Item("C").ColumnWidth = $targetColWidth
